$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.240043520927429
$ws.Range("B1").Value = 1.598049283027649
$ws.Range("C1").Value = 2.169703245162964
$ws.Range("D1").Value = 5.964247226715088
$ws.Range("E1").Value = 3.051313638687134
